# Fill in the "Lucas" / "Pronto" assignments for the
# "6.2.1.1 Diagramas de Sequência" row of the task table.
$d = $word.ActiveDocument

$table = $d.Tables.Item(1)

for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $row = $table.Rows.Item($i)
    $firstCellText = $row.Cells.Item(1).Range.Text
    if ($firstCellText -like "*Diagramas de Sequência*") {
        $row.Cells.Item(2).Range.Text = "Lucas "
        $row.Cells.Item(3).Range.Text = "Pronto"
        break
    }
}
